$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Hide gridlines on every sheet + set explicit column widths
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)   # "Table"
$ws2 = $wb.Worksheets.Item(2)   # "Variables"
$ws3 = $wb.Worksheets.Item(3)   # "Codelists"
$ws4 = $wb.Worksheets.Item(4)   # "Data"

$ws1.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws2.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws3.Activate()
$excel.ActiveWindow.DisplayGridlines = $false
$ws4.Activate()
$excel.ActiveWindow.DisplayGridlines = $false

$ws1.Activate()

# Sheet1 "Table" columns A:B
$ws1.Columns.Item(1).ColumnWidth = 15.75
$ws1.Columns.Item(2).ColumnWidth = 43.75

# Sheet2 "Variables" columns A:I
$ws2.Columns.Item(1).ColumnWidth = 7.75
$ws2.Columns.Item(2).ColumnWidth = 7.75
$ws2.Columns.Item(3).ColumnWidth = 6.75
$ws2.Columns.Item(4).ColumnWidth = 11.75
$ws2.Columns.Item(5).ColumnWidth = 16.75
$ws2.Columns.Item(6).ColumnWidth = 25.75
$ws2.Columns.Item(7).ColumnWidth = 28.75
$ws2.Columns.Item(8).ColumnWidth = 31.75
$ws2.Columns.Item(9).ColumnWidth = 31.75

# Sheet3 "Codelists" columns A:F
$ws3.Columns.Item(1).ColumnWidth = 7.75
$ws3.Columns.Item(2).ColumnWidth = 3.75
$ws3.Columns.Item(3).ColumnWidth = 8.75
$ws3.Columns.Item(4).ColumnWidth = 25.75
$ws3.Columns.Item(5).ColumnWidth = 34.75
$ws3.Columns.Item(6).ColumnWidth = 8.75

# Sheet4 "Data" columns A:D
$ws4.Columns.Item(1).ColumnWidth = 4.75
$ws4.Columns.Item(2).ColumnWidth = 5.75
$ws4.Columns.Item(3).ColumnWidth = 5.75
$ws4.Columns.Item(4).ColumnWidth = 7.75

# ---------------------------------------------------------------------------
# 2. Turn each sheet's used range into a formatted Excel Table (ListObject)
# ---------------------------------------------------------------------------
$t1 = $ws1.ListObjects.Add(1, $ws1.Range("A1:B32"), $null, 1)
$t2 = $ws2.ListObjects.Add(1, $ws2.Range("A1:I5"), $null, 1)
$t3 = $ws3.ListObjects.Add(1, $ws3.Range("A1:F11"), $null, 1)
$t4 = $ws4.ListObjects.Add(1, $ws4.Range("A1:D85"), $null, 1)

# Rename back-to-front: renaming a table to a name that collides with a
# not-yet-assigned default table name (e.g. "Table4") confuses the engine's
# auto-naming counter if done front-to-back, silently dropping an earlier
# table. Renaming from the last-added table to the first avoids any clash.
$t4.Name = "Table6"
$t3.Name = "Table5"
$t2.Name = "Table4"
$t1.Name = "Table3"
